$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1. guess_max sheet (sheet index 2): update the stale selection sqref
#    from A1:E1048576 to A1:E65536 (legacy xls row-limit artifact).
# ----------------------------------------------------------------------
$wsGuessMax = $wb.Worksheets.Item(2)
$wsGuessMax.Range("A1:E65536").Select()

# ----------------------------------------------------------------------
# 2. logical_coercion sheet (sheet index 4): rebuild/extend the data,
#    beefing up the logical-coercion tests and dropping the old,
#    low-value coercion test (the formula ="F" cell and the trailing
#    bare FALSE literal).
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item(4)

# --- Shift the still-useful rows down into their new homes, preserving
#     their original styling (notably the date cell's number format)
#     via Range.Insert() rather than rewriting them from scratch. ------

# 2a. Open a gap above the "0" row so a blank separator row lands
#     between "cabbage" and "0" in the final layout.
$ws.Range("A5").Insert()

# 2b. Open a gap above "cabbage" so it (and everything that now follows
#     it) moves down to row 14.
$ws.Range("A4:A13").Insert()

# After the two inserts:
#   row3  -> still holds the old ="F" formula (to be replaced)
#   row14 -> "cabbage"
#   row15 -> blank (the separator row created above)
#   row16 -> 0
#   row17 -> 1
#   row18 -> the date value (40908) with its original date style
#   row19 -> the old trailing FALSE literal (to be dropped)

# --- Drop the old, low-value cells -------------------------------------
$ws.Range("A19").ClearContents()

# --- Header row ----------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "logical"
$ws.Cells.Item(1, 2).Value = "explanation"

# --- Column A: logical-coercion examples, filled top to bottom ----------
$ws.Cells.Item(2, 1).Value = $true
$ws.Cells.Item(3, 1).Value = $false

$ws.Cells.Item(4, 1).Formula = "=TRUE()"
$ws.Cells.Item(5, 1).Formula = "=FALSE()"

$ws.Cells.Item(6, 1).Formula = '="true"'
$ws.Cells.Item(7, 1).Formula = '="false"'

$ws.Cells.Item(8, 1).Value = "'true"
$ws.Cells.Item(9, 1).Value = "'false"
$ws.Cells.Item(10, 1).Value = "T"
$ws.Cells.Item(11, 1).Value = "F"
$ws.Cells.Item(12, 1).Value = "'True"
$ws.Cells.Item(13, 1).Value = "'False"

# row 14 ("cabbage") and rows 16-18 (0 / 1 / date) already carried over
# correctly via the inserts above.

# --- Column B: explanatory annotations, filled top to bottom ------------
$ws.Cells.Item(2, 2).Value = "static logical"
$ws.Cells.Item(3, 2).Value = "static logical"
$ws.Cells.Item(4, 2).Value = "formula logical"
$ws.Cells.Item(5, 2).Value = "formula logical"
$ws.Cells.Item(6, 2).Value = "string logical"
$ws.Cells.Item(7, 2).Value = "string logical"
$ws.Cells.Item(8, 2).Value = "string logical"
$ws.Cells.Item(9, 2).Value = "string logical"
$ws.Cells.Item(10, 2).Value = "string logical"
$ws.Cells.Item(11, 2).Value = "string logical"
$ws.Cells.Item(12, 2).Value = "string logical"
$ws.Cells.Item(13, 2).Value = "string logical"
$ws.Cells.Item(14, 2).Value = "string not logical"
$ws.Cells.Item(15, 2).Value = "blank"
$ws.Cells.Item(16, 2).Value = "numeric"
$ws.Cells.Item(17, 2).Value = "numeric"
$ws.Cells.Item(18, 2).Value = "date"

# --- Selection / active-sheet bookkeeping -------------------------------
$ws.Activate()
$ws.Range("B18").Select()
